# This workbook's data rows (2-35, excluding 32 which is unchanged) have their
# D (Fecha), H (Variedad), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), N (Unidad de comercializacion), O (Origen)
# and P (Precio $/Kg) values redistributed among the rows (a permutation),
# while columns A, B, C, E, F, G, I, Q, R stay put.
#
# $mapping[r] = s means: row r's NEW values (for the columns above) are taken
# from row s's OLD (current, before-edit) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2 = 14
    3 = 33
    4 = 19
    5 = 25
    6 = 12
    7 = 10
    8 = 7
    9 = 21
    10 = 8
    11 = 23
    12 = 31
    13 = 30
    14 = 11
    15 = 9
    16 = 13
    17 = 3
    18 = 34
    19 = 5
    20 = 17
    21 = 20
    22 = 29
    23 = 27
    24 = 2
    25 = 4
    26 = 18
    27 = 35
    28 = 26
    29 = 24
    30 = 6
    31 = 16
    33 = 22
    34 = 28
    35 = 15
}

$cols = @("D", "H", "J", "K", "L", "M", "N", "O", "P")

# Snapshot all the current (before-edit) values for the affected columns,
# for every row, before we start overwriting anything.
$snapshot = @{}
foreach ($r in 2..35) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write the new values according to the mapping.
foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}
